$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Fill in the remaining columns for row 32 (date 12/4/2019 -> serial 43803)
$ws.Range("B32").Value = "Research YOLO in C++"
$ws.Range("C32").Value = "Started documentation + visitation company for robohub"
$ws.Range("D32").Value = "Yolo is not in c++ without wrapper so we stay with the yolo we currently use and fix the communication with Marian"

# Add two new rows with just a date value, matching the row 32 date formatting
$ws.Range("A33").Value = 43809
$ws.Range("A33").NumberFormat = $ws.Range("A32").NumberFormat

$ws.Range("A34").Value = 43810
$ws.Range("A34").NumberFormat = $ws.Range("A32").NumberFormat
